# New hind leg subassemblies
# Add "Sheet3": a Part / Mass (g) table used for center-of-mass calcs.
$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet, then rename it.
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws.Name = "Sheet3"

# --- Labels first, written in the exact order they were originally typed
#     so the workbook shared-string table is built up in the same order
#     as the authored file (Part/Mass(g) header, then each part name, the
#     "nuts (m)" header and "x" marker come later, matching the source). ---
$ws.Range("A1").Value = "Part"
$ws.Range("B1").Value = "Mass (g)"
$ws.Range("A2").Value = "BL_012v5"
$ws.Range("A25").Value = "medium nut"
$ws.Range("A3").Value = "BL_002.rat"
$ws.Range("A24").Value = "post and screw"
$ws.Range("A4").Value = "BL_003v3"
$ws.Range("A14").Value = "BO_025v6"
$ws.Range("A15").Value = "6597K8"
$ws.Range("A16").Value = "6597K19"
$ws.Range("A17").Value = "6597K118"
$ws.Range("A18").Value = "6597K119"
$ws.Range("A5").Value = "BL_011j2"
$ws.Range("A6").Value = "BL_004v4"
$ws.Range("A7").Value = "BL_005.rat"
$ws.Range("A8").Value = "BL_011j3"
$ws.Range("A9").Value = "BL_006v2"
$ws.Range("A10").Value = "BL_007v2"
$ws.Range("A11").Value = "BL_008"
$ws.Range("A12").Value = "BL_009"
$ws.Range("A13").Value = "BL_010"
$ws.Range("A26").Value = "large nut"
$ws.Range("C1").Value = "nuts (m)"
$ws.Range("E2").Value = "x"
$ws.Range("A19").Value = "encoder + mount small"
$ws.Range("A20").Value = "encoder + mount large"
$ws.Range("A21").Value = "dshaft 119mm"
$ws.Range("A22").Value = "dshaft 80mm"
$ws.Range("A23").Value = "dshaft 52mm"
$ws.Range("A27").Value = "4-40 SHS 1250"
$ws.Range("A28").Value = "4-40 SHS 0500"

# --- Part mass table: qty of nuts (col C) reduces the raw mass (col B) by
#     the nut mass constant in $B$25; col E marks rows using nuts. ---
$ws.Range("B2").Value = 44
$ws.Range("C2").Value = 2
$ws.Range("D2").Formula = "=B2-C2*`$B`$25"

$ws.Range("B3").Value = 23
$ws.Range("C3").Value = 4
$ws.Range("D3").Formula = "=B3-C3*`$B`$25"
$ws.Range("E3").Value = "x"

$ws.Range("B4").Value = 39
$ws.Range("D4").Formula = "=B4-C4*`$B`$25"
$ws.Range("E4").Value = "x"

$ws.Range("B5").Value = 13
$ws.Range("D5").Formula = "=B5-C5*`$B`$25"

$ws.Range("B6").Value = 38
$ws.Range("D6").Formula = "=B6-C6*`$B`$25"
$ws.Range("E6").Value = "x"

$ws.Range("B7").Value = 35
$ws.Range("C7").Value = 4
$ws.Range("D7").Formula = "=B7-C7*`$B`$25"
$ws.Range("E7").Value = "x"

$ws.Range("B8").Value = 7
$ws.Range("D8").Formula = "=B8-C8*`$B`$25"

$ws.Range("B9").Value = 20
$ws.Range("D9").Formula = "=B9-C9*`$B`$25"
$ws.Range("E9").Value = "x"

$ws.Range("B10").Value = 14
$ws.Range("D10").Formula = "=B10-C10*`$B`$25"
$ws.Range("E10").Value = "x"

$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 1
$ws.Range("D11").Formula = "=B11-C11*`$B`$25"
$ws.Range("E11").Value = "x"

$ws.Range("B12").Value = 8
$ws.Range("C12").Value = 1
$ws.Range("D12").Formula = "=B12-C12*`$B`$25"
$ws.Range("E12").Value = "x"

$ws.Range("B13").Value = 1
$ws.Range("D13").Formula = "=B13-C13*`$B`$25"
$ws.Range("E13").Value = "x"

$ws.Range("B14").Value = 105
$ws.Range("C14").Value = 6
$ws.Range("D14").Formula = "=B14-C14*`$B`$25"
$ws.Range("E14").Value = "x"

$ws.Range("B15").Value = 1.7
$ws.Range("D15").Formula = "=B15-C15*`$B`$25"
$ws.Range("E15").Value = "x"

$ws.Range("B16").Value = 7
$ws.Range("D16").Formula = "=B16-C16*`$B`$25"
$ws.Range("E16").Value = "x"

$ws.Range("B17").Value = 79
$ws.Range("D17").Formula = "=B17-C17*`$B`$25"
$ws.Range("E17").Value = "x"

$ws.Range("B18").Value = 79
$ws.Range("D18").Formula = "=B18-C18*`$B`$25"

$ws.Range("B19").Value = 21
$ws.Range("D19").Formula = "=B19-C19*`$B`$25"

$ws.Range("B20").Value = 28
$ws.Range("D20").Formula = "=B20-C20*`$B`$25"

$ws.Range("B21").Value = 29
$ws.Range("D21").Formula = "=B21-C21*`$B`$25"

$ws.Range("B22").Value = 20
$ws.Range("D22").Formula = "=B22-C22*`$B`$25"

$ws.Range("B23").Value = 12
$ws.Range("D23").Formula = "=B23-C23*`$B`$25"

$ws.Range("B24").Value = 2
$ws.Range("D24").Formula = "=B24-C24*`$B`$25"

# --- Reference constants (nut/screw masses, unit conversions) ---
$ws.Range("B25").Formula = "=4/10"

$ws.Range("B26").Formula = "=4/6"
$ws.Range("B26").NumberFormat = "0.0"
$ws.Range("D26").NumberFormat = "0.0"

$ws.Range("B27").Formula = "=13/10"

$ws.Range("B28").Formula = "=2/4"

# --- Widen columns A:B so the long part names are readable ---
$ws.Range("A1:B28").ColumnWidth = 26.45

# --- Sheet3 becomes the active tab / selection, matching the saved view ---
$ws.Range("F29").Select()
